# Dev IV Project Rubric.xlsx - grading entries filled in + selection moved.
#
# The author marked several more rubric rows as "Milestone I / Complete(X)"
# (columns E/F use the existing shared strings "I" and "X"), marked the
# GIT-effectiveness rows (C73:C74) as complete too, and left the cursor on
# F37 instead of F6 when they saved. The dependent SUM/SUMIF formulas in
# column G and in H4/K4/L4/H6/H8/I8/H10 recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 17, 25, 26, 27 and 28: mark as Milestone "I", Complete "X".
$rows = @(17, 25, 26, 27, 28)
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "I"
    $ws.Range("F$r").Value = "X"
}

# Effective-use-of-GIT rows: mark Milestone I column complete too.
$ws.Range("C73").Value = "X"
$ws.Range("C74").Value = "X"

# Cursor/selection left on F37 when the author saved.
[void]$ws.Range("F37").Select()
